$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.961421333333333
$ws.Range("H2").Value = 11.884264
$ws.Range("I2").Value = 0.3114993985605504
$ws.Range("J2").Value = 0.3114993985605504
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.310403
$ws.Range("N2").Value = 0.9312090000000001
$ws.Range("Q2").Value = 1.229637066130667
$ws.Range("R2").Value = 11.066733595176
$ws.Range("S2").Value = 0.3114993985605504
$ws.Range("T2").Value = 0.3114993985605504

$ws.Range("I3").Value = 0.5009735319462221
$ws.Range("J3").Value = 0.500973531946222
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.310403
$ws.Range("N3").Value = 0.9312090000000001
$ws.Range("Q3").Value = 1.977582065577334
$ws.Range("R3").Value = 17.798238590196
$ws.Range("S3").Value = 0.5009735319462221
$ws.Range("T3").Value = 0.500973531946222

$ws.Range("I4").Value = 0.1875270694932276
$ws.Range("J4").Value = 0.1875270694932276
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.310403
$ws.Range("N4").Value = 0.9312090000000001
$ws.Range("Q4").Value = 0.740259007296
$ws.Range("R4").Value = 6.662331065664
$ws.Range("S4").Value = 0.1875270694932276
$ws.Range("T4").Value = 0.1875270694932276
